$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Update the authorship template string before removing columns so the
# shared-string edit is independent of later column shifts.
# (Single-quoted so PowerShell does not treat ${...} as a variable reference.)
$ws.Range("BB2").Value = '${summary.authority}'

# Delete the taxonomic rank columns that are no longer collected:
# suborder, infraorder, superfamily (columns AR, AS, AT).
$ws.Range("AR1:AT1").EntireColumn.Delete()

# Delete the leading Taxon_Local_ID column (column A).
$ws.Range("A1").EntireColumn.Delete()
